# Regenerate the random HashMap test fixture in column A (A1:A100) with a
# new set of random values (1-10), per "add test for hashN method".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    2,1,2,7,3,1,8,6,2,3,
    6,8,9,10,8,9,9,9,5,5,
    6,9,10,5,1,4,2,9,2,5,
    3,6,5,6,3,2,8,7,5,8,
    1,10,10,5,3,10,4,3,10,9,
    2,3,8,7,6,2,2,5,1,1,
    5,8,3,1,6,5,4,10,5,7,
    4,10,4,1,6,5,2,7,6,9,
    4,9,7,9,3,8,1,8,2,10,
    9,6,4,9,3,1,5,10,4,5
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}
